# Respawn cover letter edit:
#  1. Remove the placeholder paragraph "Talk about wanting to create power
#     fantasy that they go for here? -".
#  2. Merge the "Combat, weapons, action..." paragraph with the following
#     "The free-flowing movement of Titanfall..." paragraph (drop the
#     paragraph break between them).
#  3. Rewrite the final "Games with..." paragraph with the new, fleshed-out
#     studio-praise text, tighten its spacing to 9pt after, and move the
#     hidden "_GoBack" bookmark from its old spot (next to "semester") to
#     its new spot (inside "fan|tasies!").

$d = $word.ActiveDocument

# --- Step 1: delete the "Talk about..." placeholder paragraph -------------
$talkPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Talk about wanting to create power fantasy")) {
        $talkPara = $d.Paragraphs($i)
        break
    }
}
$talkPara.Range.Delete()

# --- Step 2: merge the Combat / Titanfall paragraphs ----------------------
$combatPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Combat, weapons, action and animation")) {
        $combatPara = $d.Paragraphs($i)
        break
    }
}
$paraMark = $d.Range($combatPara.Range.End - 1, $combatPara.Range.End)
$paraMark.Delete()

# --- Step 3: rewrite the closing "Games with..." paragraph ----------------
$closingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Games with fluid, fast-paced action")) {
        $closingPara = $d.Paragraphs($i)
        break
    }
}
$closingContent = $d.Range($closingPara.Range.Start, $closingPara.Range.End - 1)
$closingContent.Delete()

$rsquo = [char]0x2019
$newTextBeforeBookmark = "Games with fluid, fast-paced action are what Respawn is the flag-bearer of! " + `
    "I like the studio" + $rsquo + "s focus on pushing the envelope for games in new directions and the emphasis placed on gameplay feel. " + `
    "I admire how you have created such innovative gameplay that is fluid, responsive and achieved a very cool power fantasy with it too! " + `
    "It is the kind of player experience that I would like to be creating myself. " + `
    "And I would absolutely love to be a part of Respawn Entertainment and play my part in creating games that are trailblazers of these creative fan"
$newTextAfterBookmark = "tasies!"

$closingPara = $d.Paragraphs($closingPara.Index)
$insertStart = $closingPara.Range.Start
$closingPara.Range.InsertBefore($newTextBeforeBookmark + $newTextAfterBookmark)

$bookmarkPos = $insertStart + $newTextBeforeBookmark.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$closingPara = $d.Paragraphs($closingPara.Index)
$closingPara.SpaceAfter = 9
